# Revision: add the "chemical_recycling_pyrolysis" parameter row right
# after the existing "chemical_recycling_gasification" row (row 9),
# pushing every row below it down by one — matching the diff, which
# shows a new <row r="10"> inserted and all subsequent rows (old 10-24,
# now 11-25) shifted down by one, with the dimension growing from
# A1:C24 to A1:C25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10 (everything from row 10 down,
# including "fossil_routes", shifts to row 11, etc.)
$ws.Rows("10:10").Insert()

# Populate the newly inserted row with the new parameter.
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
